# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" for ac0b8092-787e-4c1c-9179-d816e2c0177f.md
# on the zh-cn and de-de sheets, and the corresponding "Latest HO Xliff Generate Date"
# on the Overview sheet (set to the max/most-recent of the two locale handoff times).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 7 is ac0b8092-787e-4c1c-9179-d816e2c0177f.md, column H = Latest Handoff Datetime
$wsZhCn.Range("H7").Value = "2016-08-21 14:49:26"

# de-de sheet: row 7 is ac0b8092-787e-4c1c-9179-d816e2c0177f.md, column H = Latest Handoff Datetime
$wsDeDe.Range("H7").Value = "2016-08-21 14:49:30"

# Overview sheet: row 7 is ac0b8092-787e-4c1c-9179-d816e2c0177f.md, column G = Latest HO Xliff Generate Date
$wsOverview.Range("G7").Value = "2016-08-21 14:49:30"
